$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Add two new rows to the "Gemeinde/Vorhaben/Standort" table (Table 2):
#    "Einsprache" (dropdown "n") + EINSPRECHENDE for-loop
#    "Rechtsverwahrung-" (dropdown "en") + RECHTSVERWAHRENDE for-loop
# ---------------------------------------------------------------------------
$t = $d.Tables(2)

$row1 = $t.Rows.Add()
$row1Cell1Xml = @'
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr><w:t>Einsprache</w:t></w:r><w:r><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr><w:fldChar w:fldCharType="begin"><w:ffData><w:name w:val=""/><w:enabled/><w:calcOnExit w:val="0"/><w:ddList><w:listEntry w:val="n"/></w:ddList></w:ffData></w:fldChar></w:r><w:r><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr><w:instrText xml:space="preserve"> FORMDROPDOWN </w:instrText></w:r><w:r><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr></w:r><w:r><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p></w:body></w:document>
'@
$row1.Cells(1).Range.Paragraphs(1).Range.InsertXML($row1Cell1Xml)

$row1Cell2Xml = @'
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="AufzhlungBrief"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr><w:ind w:left="357" w:hanging="357"/><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr><w:t xml:space="preserve">{% </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr><w:t>for</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr><w:t xml:space="preserve"> POSITION in EINSPRECHENDE %}</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="AufzhlungBrief"/></w:pPr><w:r><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr><w:t xml:space="preserve">{{POSITION.NAME}}, {{POSITION.ADRESSE}}{% </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr><w:t>endfor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr><w:t xml:space="preserve"> %}</w:t></w:r></w:p></w:body></w:document>
'@
$row1.Cells(2).Range.Paragraphs(1).Range.InsertXML($row1Cell2Xml)

$row2 = $t.Rows.Add()
$row2Cell1Xml = @'
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr><w:t>Rechtsverwahrung-</w:t></w:r><w:r><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr><w:fldChar w:fldCharType="begin"><w:ffData><w:name w:val=""/><w:enabled/><w:calcOnExit w:val="0"/><w:ddList><w:listEntry w:val="en"/></w:ddList></w:ffData></w:fldChar></w:r><w:r><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr><w:instrText xml:space="preserve"> FORMDROPDOWN </w:instrText></w:r><w:r><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr></w:r><w:r><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p></w:body></w:document>
'@
$row2.Cells(1).Range.Paragraphs(1).Range.InsertXML($row2Cell1Xml)

$row2Cell2Xml = @'
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="AufzhlungBrief"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr><w:ind w:left="357" w:hanging="357"/><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr><w:t xml:space="preserve">{% </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr><w:t>for</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr><w:t xml:space="preserve"> POSITION in RECHTSVERWAHRENDE %}</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="AufzhlungBrief"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr><w:ind w:left="357" w:hanging="357"/><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr><w:t xml:space="preserve">{{POSITION.NAME}}, {{POSITION.ADRESSE}}{% </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr><w:t>endfor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr><w:t xml:space="preserve"> %}</w:t></w:r></w:p></w:body></w:document>
'@
$row2.Cells(2).Range.Paragraphs(1).Range.InsertXML($row2Cell2Xml)

Write-Host "Step 1 done: table rows added"
